$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell "Save" in H1, matching the style of the existing header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the "Save" column values for rows 2-11
$saveValues = @(0, 1, 0, 1, 0, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
